$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.459.20'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '1.948.84'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '243.30'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.625'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '59.75'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.05%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.378'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +4.48%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0788'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -8.04%  '
$ws.Range("E11").Value = '  +0.46%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '14.13'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.63%  '
$ws.Range("D13").Value = '2.235.06'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.822'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.27%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.32'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.10%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.40%  '
$ws.Range("D17").Value = '1.943.65'
$ws.Range("E17").Value = '  +0.48%  '
$ws.Range("D18").Value = '36.395.37'
$ws.Range("E18").Value = '  +0.01%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.22%  '
$ws.Range("D20").Value = '0.0₃0847'
$ws.Range("E20").Value = '  -1.74%  '
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '228.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '
$ws.Range("B22").Value = 'Uniswap'
$ws.Range("C22").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.29%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.43'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.94%  '
$ws.Range("E25").Value = '  +2.66%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.142'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +9.04%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.10'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.81'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.19'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.30'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +18.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.120'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.48%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.73'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0610'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.42'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.91%  '
$ws.Range("B35").Value = 'BinanceUSD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.05%  '
$ws.Range("B36").Value = 'LidoDAOToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.42'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.75%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.77'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.40'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -12.41%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0966'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.28%  '
$ws.Range("E41").Value = '  -0.82%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.16'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.82%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0209'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.37%  '
$ws.Range("D45").Value = '1.359.54'
$ws.Range("E45").Value = '  +1.74%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '88.28'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.68%  '
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.10'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.82'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '44.94'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.07%  '
$ws.Range("D51").Value = '2.130.57'
$ws.Range("E51").Value = '  +0.68%  '
